$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5616.3335
$ws.Range("I40").Value = 6949.75
$ws.Range("K40").Value = 6949.75
$ws.Range("M40").Value = -6774.75
$ws.Range("H43").Value = 11145
$ws.Range("J43").Value = 3683.9
$ws.Range("L43").Value = 3683.9
$ws.Range("N43").Value = -3821.9
$ws.Range("H92").Value = 247.36363
$ws.Range("I92").Value = 247.36363
$ws.Range("K92").Value = 247.36363
$ws.Range("M92").Value = 1000.63637
$ws.Range("H101").Value = 220
$ws.Range("I101").Value = 204
$ws.Range("J101").Value = 300
$ws.Range("K101").Value = 612
$ws.Range("L101").Value = 900
$ws.Range("M101").Value = 1010
$ws.Range("N101").Value = -4144
$ws.Range("H111").Value = 2085.4
$ws.Range("I111").Value = 816.3333
$ws.Range("K111").Value = 2448.9999
$ws.Range("M111").Value = 618.0001000000002
$ws.Range("H116").Value = 29208.5
$ws.Range("I116").Value = 5161.6665
$ws.Range("K116").Value = 5161.6665
$ws.Range("M116").Value = -1719.6665

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 228187.62
$ws.Range("I32").Value = 233080.27
$ws.Range("K32").Value = 233080.27
$ws.Range("M32").Value = -232793.27
$ws.Range("H61").Value = 16250
$ws.Range("I61").Value = 12500
$ws.Range("K61").Value = 12500
$ws.Range("M61").Value = -12288
$ws.Range("H63").Value = 2133.7368
$ws.Range("I63").Value = 2155.5293
$ws.Range("J63").Value = 1948.5
$ws.Range("K63").Value = 2155.5293
$ws.Range("L63").Value = 1948.5
$ws.Range("M63").Value = -1469.5293
$ws.Range("N63").Value = -3320.5
$ws.Range("H66").Value = 2133.7368
$ws.Range("I66").Value = 2155.5293
$ws.Range("J66").Value = 1948.5
$ws.Range("K66").Value = 10777.6465
$ws.Range("L66").Value = 9742.5
$ws.Range("M66").Value = -7345.646500000001
$ws.Range("N66").Value = -16606.5
$ws.Range("H74").Value = 5395.3193
$ws.Range("I74").Value = 817.8611
$ws.Range("K74").Value = 817.8611
$ws.Range("M74").Value = 56.13890000000004
$ws.Range("H77").Value = 5395.3193
$ws.Range("I77").Value = 817.8611
$ws.Range("K77").Value = 4089.3055
$ws.Range("M77").Value = 278.6945000000001
$ws.Range("H122").Value = 2909.9333
$ws.Range("I122").Value = 2832.3635
$ws.Range("J122").Value = 3123.25
$ws.Range("K122").Value = 8497.0905
$ws.Range("L122").Value = 9369.75
$ws.Range("M122").Value = -6047.0905
$ws.Range("N122").Value = -14269.75
$ws.Range("H136").Value = 16250
$ws.Range("I136").Value = 12500
$ws.Range("K136").Value = 37500
$ws.Range("M136").Value = -34950

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 153.3
$ws.Range("I7").Value = 113.2
$ws.Range("J7").Value = 193.4
$ws.Range("K7").Value = 113.2
$ws.Range("L7").Value = 193.4
$ws.Range("M7").Value = -0.2000000000000028
$ws.Range("N7").Value = -419.4
$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 1999
$ws.Range("K62").Value = 1999
$ws.Range("M62").Value = -1375
$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 1999
$ws.Range("K65").Value = 9995
$ws.Range("M65").Value = -6875
$ws.Range("H92").Value = 100000
$ws.Range("J92").Value = 100000
$ws.Range("L92").Value = 100000
$ws.Range("N92").Value = -104992
$ws.Range("H105").Value = 12647.111
$ws.Range("I105").Value = 14989.143
$ws.Range("K105").Value = 14989.143
$ws.Range("M105").Value = -13242.143

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 205.5
$ws.Range("J23").Value = 205.5
$ws.Range("L23").Value = 616.5
$ws.Range("N23").Value = -1086.5
$ws.Range("H68").Value = 6313.1816
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 6313.1816
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 18939.5448
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -20561.5448
$ws.Range("H71").Value = 6313.1816
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 6313.1816
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 56818.6344
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -64930.6344
$ws.Range("H80").Value = 2830.6667
$ws.Range("J80").Value = 3250
$ws.Range("L80").Value = 9750
$ws.Range("N80").Value = -11622
$ws.Range("H83").Value = 2830.6667
$ws.Range("J83").Value = 3250
$ws.Range("L83").Value = 29250
$ws.Range("N83").Value = -38610
$ws.Range("H92").Value = 660
$ws.Range("I92").Value = 580
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 1740
$ws.Range("L92").Value = 2700
$ws.Range("M92").Value = -492
$ws.Range("N92").Value = -5196
$ws.Range("H121").Value = 2987.375
$ws.Range("I121").Value = 364.375
$ws.Range("K121").Value = 1093.125
$ws.Range("M121").Value = 216.875
$ws.Range("H129").Value = 1253928.2
$ws.Range("J129").Value = 4750
$ws.Range("L129").Value = 14250
$ws.Range("N129").Value = -24250
$ws.Range("H132").Value = 748
$ws.Range("I132").Value = 546.875
$ws.Range("K132").Value = 4921.875
$ws.Range("M132").Value = -2391.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2309.8
$ws.Range("J22").Value = 2309.8
$ws.Range("L22").Value = 2309.8
$ws.Range("N22").Value = -2899.8
$ws.Range("H27").Value = 2309.8
$ws.Range("J27").Value = 2309.8
$ws.Range("L27").Value = 2309.8
$ws.Range("N27").Value = -2523.8
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H82").Value = 3186
$ws.Range("I82").Value = 4279.75
$ws.Range("J82").Value = 998.5
$ws.Range("K82").Value = 4279.75
$ws.Range("L82").Value = 998.5
$ws.Range("M82").Value = -3918.75
$ws.Range("N82").Value = -1720.5
$ws.Range("H85").Value = 3186
$ws.Range("I85").Value = 4279.75
$ws.Range("J85").Value = 998.5
$ws.Range("K85").Value = 4279.75
$ws.Range("L85").Value = 998.5
$ws.Range("M85").Value = -3031.75
$ws.Range("N85").Value = -3494.5
$ws.Range("H125").Value = 75715
$ws.Range("J125").Value = 75715
$ws.Range("L125").Value = 75715
$ws.Range("N125").Value = -85555
$ws.Range("H132").Value = 1390865.1
$ws.Range("I132").Value = 1787643.8
$ws.Range("J132").Value = 2140.25
$ws.Range("K132").Value = 5362931.4
$ws.Range("L132").Value = 6420.75
$ws.Range("M132").Value = -5360401.4
$ws.Range("N132").Value = -11480.75
$ws.Range("H136").Value = 9162.736999999999
$ws.Range("I136").Value = 4344.7
$ws.Range("K136").Value = 13034.1
$ws.Range("M136").Value = -10484.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1001700.6
$ws.Range("J3").Value = 2125.75
$ws.Range("L3").Value = 2125.75
$ws.Range("N3").Value = -2353.75
$ws.Range("H19").Value = 652.5
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = 1300
$ws.Range("K19").Value = 5
$ws.Range("L19").Value = 1300
$ws.Range("M19").Value = 169
$ws.Range("N19").Value = -1648
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H62").Value = 4427.6665
$ws.Range("I62").Value = 1789
$ws.Range("K62").Value = 1789
$ws.Range("M62").Value = -1165
$ws.Range("H65").Value = 4427.6665
$ws.Range("I65").Value = 1789
$ws.Range("K65").Value = 8945
$ws.Range("M65").Value = -5825
$ws.Range("H122").Value = 36170.332
$ws.Range("I122").Value = 2000.75
$ws.Range("K122").Value = 6002.25
$ws.Range("M122").Value = -3552.25
$ws.Range("H126").Value = 2330.7778
$ws.Range("I126").Value = 2389.6428
$ws.Range("K126").Value = 7168.928400000001
$ws.Range("M126").Value = -4698.928400000001
$ws.Range("H136").Value = 1608.4584
$ws.Range("I136").Value = 1589.1111
$ws.Range("K136").Value = 4767.3333
$ws.Range("M136").Value = -2217.3333
